# Update "想去人数" (want-to-go count) figures for several events.
# Sheet "展览" (Exhibitions) and sheet "全部类型" (All types) both list the
# same events, so each change must be applied in both places.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsAll  = $wb.Worksheets.Item("全部类型")

# Row numbers in "展览" and their counterparts in "全部类型"
$wsExpo.Range("F2").Value  = 261   # 蜀山·银泰百货高新店-2024漫趣地带嘉年华（免费）
$wsExpo.Range("F4").Value  = 13    # 合肥·星域动漫游戏嘉年华
$wsExpo.Range("F5").Value  = 6623  # 合肥·第十五届次元之门动漫游戏博览会
$wsExpo.Range("F6").Value  = 5395  # 合肥·首届AT次元时代动漫游戏嘉年华
$wsExpo.Range("F7").Value  = 448   # 合肥·Holic动漫游戏展
$wsExpo.Range("F8").Value  = 66    # 合肥·乐帮•崩原铁绝only同人首展
$wsExpo.Range("F12").Value = 24    # 合肥·第九届环形宇宙动漫游戏嘉年华

$wsAll.Range("F2").Value  = 261    # 蜀山·银泰百货高新店-2024漫趣地带嘉年华（免费）
$wsAll.Range("F4").Value  = 13     # 合肥·星域动漫游戏嘉年华
$wsAll.Range("F5").Value  = 6623   # 合肥·第十五届次元之门动漫游戏博览会
$wsAll.Range("F6").Value  = 5395   # 合肥·首届AT次元时代动漫游戏嘉年华
$wsAll.Range("F7").Value  = 448    # 合肥·Holic动漫游戏展
$wsAll.Range("F8").Value  = 66     # 合肥·乐帮•崩原铁绝only同人首展
$wsAll.Range("F14").Value = 24     # 合肥·第九届环形宇宙动漫游戏嘉年华
